$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended after the existing 68 rows (11.6.18 dataset).
# Columns: A=Label, B=consensus, C=PB1, D=PB2, E=PB3, F=PB4, G=PB5, H=PB6, I=PB7, J=consensus2
$newRows = @(
    @("11.6.18.1.1", "no flip"),
    @("11.6.18.1.2", "no flip"),
    @("11.6.18.2.1", "flip AP"),
    @("11.6.18.2.2", "flip DV"),
    @("11.6.18.2.3", "flip DV"),
    @("11.6.18.3.1", "no flip"),
    @("11.6.18.3.2", "no flip"),
    @("11.6.18.3.3", "no flip")
)

$startRow = 69
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $label = $newRows[$i][0]
    $consensus = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $label
    $ws.Cells.Item($r, 2).Value = $consensus
    $ws.Cells.Item($r, 3).Value = $consensus
    $ws.Cells.Item($r, 4).Value = $consensus
    $ws.Cells.Item($r, 5).Value = $consensus
    $ws.Cells.Item($r, 6).Value = "empty"
    $ws.Cells.Item($r, 7).Value = "empty"
    $ws.Cells.Item($r, 8).Value = "empty"
    $ws.Cells.Item($r, 9).Value = "empty"
    $ws.Cells.Item($r, 10).Value = $consensus
}

# Match the saved viewport/selection state from the edit session.
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("N70").Select()
